$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D (Fecha), H (Variedad), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# N (Unidad de comercializacion), P (Precio $/Kg)
$rows = @{
    2  = @{ D = 44468; H = "Verde";            J = 500; K = 1800; L = 2000; M = 1920; N = "`$/kilo";    P = 1920 }
    3  = @{ D = 44519; H = "Sin especificar";  J = 250; K = 1200; L = 1300; M = 1240; N = "`$/kilo";    P = 1240 }
    5  = @{ D = 44510; H = "Sin especificar";  J = 600; K = 1300; L = 1400; M = 1350; N = "`$/kilo";    P = 1350 }
    6  = @{ D = 44477; H = "Sin especificar";  J = 500; K = 1400; L = 1500; M = 1460; N = "`$/kilo";    P = 1460 }
    7  = @{ D = 44526; H = "Sin especificar";  J = 100; K = 1500; L = 1600; M = 1550; N = "`$/kilo";    P = 1550 }
    8  = @{ D = 44489; H = "Sin especificar";  J = 600; K = 1400; L = 1500; M = 1450; N = "`$/kilo";    P = 1450 }
    9  = @{ D = 44496; H = "Sin especificar";  J = 550; K = 1500; L = 2000; M = 1773; N = "`$/paquete"; P = 1773 }
    10 = @{ D = 44545; H = "Sin especificar";  J = 550; K = 1700; L = 1800; M = 1755; N = "`$/kilo";    P = 1755 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($r, 8).Value  = $vals.H   # H - Variedad
    $ws.Cells.Item($r, 10).Value = $vals.J   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $vals.K   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals.L   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals.M   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $vals.N   # N - Unidad de comercializacion
    $ws.Cells.Item($r, 16).Value = $vals.P   # P - Precio $/Kg
}

$wb.Save()
